$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "PnucB683"
$ws.Range("B2").Value = 23110202
$ws.Range("C2").Value = "iagevts83"
$ws.Range("D2").Value = "k#9F!U5p"
$ws.Range("F2").Value = "JsKDMsBN"
$ws.Range("G2").Value = "tFYn"
